$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 2.629231666666667
$ws.Range("H2").Value2 = 7.887695
$ws.Range("I2").Value2 = 0.1414315557047068
$ws.Range("J2").Value2 = 0.1414315557047067
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.02536566666666666
$ws.Range("N2").Value2 = 0.076097
$ws.Range("O2").Value2 = 0.04869294909818511
$ws.Range("P2").Value2 = 0.04869294909818511
$ws.Range("Q2").Value2 = 0.0666922140461111
$ws.Range("R2").Value2 = 0.600229926415
$ws.Range("S2").Value2 = 0.006886719542806417
$ws.Range("T2").Value2 = 0.006886719542806416
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 2.629231666666667
$ws.Range("H3").Value2 = 7.887695
$ws.Range("I3").Value2 = 0.1414315557047068
$ws.Range("J3").Value2 = 0.1414315557047067
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.4080956666666666
$ws.Range("N3").Value2 = 1.224287
$ws.Range("O3").Value2 = 0.783396777436295
$ws.Range("P3").Value2 = 0.783396777436295
$ws.Range("Q3").Value2 = 1.072978049829444
$ws.Range("R3").Value2 = 9.656802448465
$ws.Range("S3").Value2 = 0.1107970249668691
$ws.Range("T3").Value2 = 0.1107970249668691
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 2.629231666666667
$ws.Range("H4").Value2 = 7.887695
$ws.Range("I4").Value2 = 0.1414315557047068
$ws.Range("J4").Value2 = 0.1414315557047067
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.08746966666666667
$ws.Range("N4").Value2 = 0.262409
$ws.Range("O4").Value2 = 0.1679102734655197
$ws.Range("P4").Value2 = 0.1679102734655197
$ws.Range("Q4").Value2 = 0.2299780174727778
$ws.Range("R4").Value2 = 2.069802157255
$ws.Range("S4").Value2 = 0.0237478111950312
$ws.Range("T4").Value2 = 0.02374781119503119
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 11.42765333333333
$ws.Range("H5").Value2 = 34.28296
$ws.Range("I5").Value2 = 0.6147160060020365
$ws.Range("J5").Value2 = 0.6147160060020365
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.02536566666666666
$ws.Range("N5").Value2 = 0.076097
$ws.Range("O5").Value2 = 0.04869294909818511
$ws.Range("P5").Value2 = 0.04869294909818511
$ws.Range("Q5").Value2 = 0.2898700452355555
$ws.Range("R5").Value2 = 2.60883040712
$ws.Range("S5").Value2 = 0.02993233519009681
$ws.Range("T5").Value2 = 0.02993233519009681
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 11.42765333333333
$ws.Range("H6").Value2 = 34.28296
$ws.Range("I6").Value2 = 0.6147160060020365
$ws.Range("J6").Value2 = 0.6147160060020365
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.4080956666666666
$ws.Range("N6").Value2 = 1.224287
$ws.Range("O6").Value2 = 0.783396777436295
$ws.Range("P6").Value2 = 0.783396777436295
$ws.Range("Q6").Value2 = 4.663575805502222
$ws.Range("R6").Value2 = 41.97218224952
$ws.Range("S6").Value2 = 0.4815665381405055
$ws.Range("T6").Value2 = 0.4815665381405055
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 11.42765333333333
$ws.Range("H7").Value2 = 34.28296
$ws.Range("I7").Value2 = 0.6147160060020365
$ws.Range("J7").Value2 = 0.6147160060020365
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.08746966666666667
$ws.Range("N7").Value2 = 0.262409
$ws.Range("O7").Value2 = 0.1679102734655197
$ws.Range("P7").Value2 = 0.1679102734655197
$ws.Range("Q7").Value2 = 0.9995730278488889
$ws.Range("R7").Value2 = 8.996157250640001
$ws.Range("S7").Value2 = 0.103217132671434
$ws.Range("T7").Value2 = 0.103217132671434
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 4.24731
$ws.Range("H8").Value2 = 12.74193
$ws.Range("I8").Value2 = 0.2284711798035388
$ws.Range("J8").Value2 = 0.2284711798035388
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 0.02536566666666666
$ws.Range("N8").Value2 = 0.076097
$ws.Range("O8").Value2 = 0.04869294909818511
$ws.Range("P8").Value2 = 0.04869294909818511
$ws.Range("Q8").Value2 = 0.10773584969
$ws.Range("R8").Value2 = 0.96962264721
$ws.Range("S8").Value2 = 0.01112493552857601
$ws.Range("T8").Value2 = 0.01112493552857601
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 4.24731
$ws.Range("H9").Value2 = 12.74193
$ws.Range("I9").Value2 = 0.2284711798035388
$ws.Range("J9").Value2 = 0.2284711798035388
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.4080956666666666
$ws.Range("N9").Value2 = 1.224287
$ws.Range("O9").Value2 = 0.783396777436295
$ws.Range("P9").Value2 = 0.783396777436295
$ws.Range("Q9").Value2 = 1.73330880599
$ws.Range("R9").Value2 = 15.59977925391
$ws.Range("S9").Value2 = 0.1789835859951606
$ws.Range("T9").Value2 = 0.1789835859951606
$ws.Range("D10").Value2 = "MuSCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 4.24731
$ws.Range("H10").Value2 = 12.74193
$ws.Range("I10").Value2 = 0.2284711798035388
$ws.Range("J10").Value2 = 0.2284711798035388
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.08746966666666667
$ws.Range("N10").Value2 = 0.262409
$ws.Range("O10").Value2 = 0.1679102734655197
$ws.Range("P10").Value2 = 0.1679102734655197
$ws.Range("Q10").Value2 = 0.37151078993
$ws.Range("R10").Value2 = 3.34359710937
$ws.Range("S10").Value2 = 0.03836265827980213
$ws.Range("T10").Value2 = 0.03836265827980213
$ws.Range("D11").Value2 = "ECs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 0.2859396666666667
$ws.Range("H11").Value2 = 0.8578190000000001
$ws.Range("I11").Value2 = 0.01538125848971795
$ws.Range("J11").Value2 = 0.01538125848971795
$ws.Range("K11").Value2 = 2
$ws.Range("L11").Value2 = 0.6666666666666666
$ws.Range("M11").Value2 = 0.02536566666666666
$ws.Range("N11").Value2 = 0.076097
$ws.Range("O11").Value2 = 0.04869294909818511
$ws.Range("P11").Value2 = 0.04869294909818511
$ws.Range("Q11").Value2 = 0.007253050271444445
$ws.Range("R11").Value2 = 0.06527745244300001
$ws.Range("S11").Value2 = 0.0007489588367058638
$ws.Range("T11").Value2 = 0.0007489588367058638
$ws.Range("D12").Value2 = "FAPs"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 0.2859396666666667
$ws.Range("H12").Value2 = 0.8578190000000001
$ws.Range("I12").Value2 = 0.01538125848971795
$ws.Range("J12").Value2 = 0.01538125848971795
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 0.4080956666666666
$ws.Range("N12").Value2 = 1.224287
$ws.Range("O12").Value2 = 0.783396777436295
$ws.Range("P12").Value2 = 0.783396777436295
$ws.Range("Q12").Value2 = 0.1166907388947778
$ws.Range("R12").Value2 = 1.050216650053
$ws.Range("S12").Value2 = 0.0120496283337597
$ws.Range("T12").Value2 = 0.0120496283337597
$ws.Range("D13").Value2 = "MuSCs"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 0.2859396666666667
$ws.Range("H13").Value2 = 0.8578190000000001
$ws.Range("I13").Value2 = 0.01538125848971795
$ws.Range("J13").Value2 = 0.01538125848971795
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 0.08746966666666667
$ws.Range("N13").Value2 = 0.262409
$ws.Range("O13").Value2 = 0.1679102734655197
$ws.Range("P13").Value2 = 0.1679102734655197
$ws.Range("Q13").Value2 = 0.02501104733011111
$ws.Range("R13").Value2 = 0.225099425971
$ws.Range("S13").Value2 = 0.002582671319252388
$ws.Range("T13").Value2 = 0.002582671319252388
